# Auto-generated update of Dades_Meteo sheet (resum diari meteocat)
# Commit: Update automàtic: dades i banners [2026-02-11 23:20]
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

# Row 2
$ws.Cells.Item(2, 5).Value = "2026-02-11 23:18:31"

# Row 3
$ws.Cells.Item(3, 5).Value = "2026-02-11 23:18:33"

# Row 4
$ws.Cells.Item(4, 5).Value = "2026-02-11 23:18:36"
$ws.Cells.Item(4, 8).Value = "'59%"
$ws.Cells.Item(4, 10).Value = "1001.6 hPa"

# Row 5
$ws.Cells.Item(5, 5).Value = "2026-02-11 23:18:39"

# Row 6
$ws.Cells.Item(6, 5).Value = "2026-02-11 23:18:41"
$ws.Cells.Item(6, 10).Value = "1002.1 hPa"

# Row 7
$ws.Cells.Item(7, 5).Value = "2026-02-11 23:18:44"
$ws.Cells.Item(7, 10).Value = "1003.0 hPa"
$ws.Cells.Item(7, 15).Value = "18.5 °C"

# Row 8
$ws.Cells.Item(8, 5).Value = "2026-02-11 23:18:46"
$ws.Cells.Item(8, 8).Value = "'61%"
$ws.Cells.Item(8, 10).Value = "1002.7 hPa"

# Row 9
$ws.Cells.Item(9, 5).Value = "2026-02-11 23:18:49"
$ws.Cells.Item(9, 15).Value = "12.0 °C"

# Row 10
$ws.Cells.Item(10, 5).Value = "2026-02-11 23:18:51"

# Row 11
$ws.Cells.Item(11, 5).Value = "2026-02-11 23:18:54"
$ws.Cells.Item(11, 8).Value = "'83%"

# Row 12
$ws.Cells.Item(12, 5).Value = "2026-02-11 23:18:56"

# Row 13
$ws.Cells.Item(13, 5).Value = "2026-02-11 23:18:59"
$ws.Cells.Item(13, 7).Value = "2 cm"
$ws.Cells.Item(13, 10).Value = "1004.6 hPa"
$ws.Cells.Item(13, 15).Value = "7.1 °C"

# Row 14
$ws.Cells.Item(14, 5).Value = "2026-02-11 23:19:01"
$ws.Cells.Item(14, 8).Value = "'55%"
$ws.Cells.Item(14, 14).Value = "11.1 °C 22:32 TU"
$ws.Cells.Item(14, 15).Value = "18.1 °C"

# Row 15
$ws.Cells.Item(15, 5).Value = "2026-02-11 23:19:04"

# Row 16
$ws.Cells.Item(16, 5).Value = "2026-02-11 23:19:06"
$ws.Cells.Item(16, 15).Value = "-0.5 °C"

# Row 17
$ws.Cells.Item(17, 5).Value = "2026-02-11 23:19:09"
$ws.Cells.Item(17, 9).Value = "7.3 mm"

# Row 18
$ws.Cells.Item(18, 5).Value = "2026-02-11 23:19:11"
$ws.Cells.Item(18, 10).Value = "1002.2 hPa"

# Row 19
$ws.Cells.Item(19, 5).Value = "2026-02-11 23:19:14"
$ws.Cells.Item(19, 8).Value = "'82%"
$ws.Cells.Item(19, 15).Value = "8.9 °C"

# Row 20
$ws.Cells.Item(20, 5).Value = "2026-02-11 23:19:16"
$ws.Cells.Item(20, 12).Value = "99.0 km/h - 249º 22:43 TU"

# Row 21
$ws.Cells.Item(21, 5).Value = "2026-02-11 23:19:19"
$ws.Cells.Item(21, 10).Value = "1005.0 hPa"
$ws.Cells.Item(21, 15).Value = "8.1 °C"

# Row 22
$ws.Cells.Item(22, 5).Value = "2026-02-11 23:19:21"
$ws.Cells.Item(22, 8).Value = "'93%"
$ws.Cells.Item(22, 9).Value = "3.9 mm"

# Row 23
$ws.Cells.Item(23, 5).Value = "2026-02-11 23:19:24"

# Row 24
$ws.Cells.Item(24, 5).Value = "2026-02-11 23:19:26"
$ws.Cells.Item(24, 8).Value = "'79%"
$ws.Cells.Item(24, 10).Value = "1006.1 hPa"

# Row 25
$ws.Cells.Item(25, 5).Value = "2026-02-11 23:19:29"
$ws.Cells.Item(25, 14).Value = "-0.6 °C 22:59 TU"
$ws.Cells.Item(25, 15).Value = "1.5 °C"

# Row 26
$ws.Cells.Item(26, 5).Value = "2026-02-11 23:19:32"
$ws.Cells.Item(26, 8).Value = "'69%"
$ws.Cells.Item(26, 10).Value = "1002.0 hPa"

# Row 27
$ws.Cells.Item(27, 5).Value = "2026-02-11 23:19:34"

# Row 28
$ws.Cells.Item(28, 5).Value = "2026-02-11 23:19:37"
$ws.Cells.Item(28, 10).Value = "1002.3 hPa"
$ws.Cells.Item(28, 15).Value = "11.2 °C"

# Row 29
$ws.Cells.Item(29, 5).Value = "2026-02-11 23:19:39"
$ws.Cells.Item(29, 15).Value = "13.4 °C"

# Row 30
$ws.Cells.Item(30, 5).Value = "2026-02-11 23:19:42"
$ws.Cells.Item(30, 8).Value = "'89%"
$ws.Cells.Item(30, 10).Value = "1002.3 hPa"
$ws.Cells.Item(30, 15).Value = "11.7 °C"

# Row 31
$ws.Cells.Item(31, 5).Value = "2026-02-11 23:19:44"
$ws.Cells.Item(31, 8).Value = "'67%"
$ws.Cells.Item(31, 10).Value = "1001.5 hPa"

# Row 32
$ws.Cells.Item(32, 5).Value = "2026-02-11 23:19:47"

# Row 33
$ws.Cells.Item(33, 5).Value = "2026-02-11 23:19:50"
$ws.Cells.Item(33, 10).Value = "1004.2 hPa"

# Row 34
$ws.Cells.Item(34, 5).Value = "2026-02-11 23:19:52"

# Row 35
$ws.Cells.Item(35, 5).Value = "2026-02-11 23:19:54"

# Row 36
$ws.Cells.Item(36, 5).Value = "2026-02-11 23:19:57"
$ws.Cells.Item(36, 10).Value = "1002.4 hPa"
$ws.Cells.Item(36, 15).Value = "13.0 °C"

# Row 37
$ws.Cells.Item(37, 5).Value = "2026-02-11 23:20:00"
$ws.Cells.Item(37, 10).Value = "1003.7 hPa"
$ws.Cells.Item(37, 15).Value = "9.2 °C"

# Row 38
$ws.Cells.Item(38, 5).Value = "2026-02-11 23:20:02"
$ws.Cells.Item(38, 8).Value = "'63%"
$ws.Cells.Item(38, 15).Value = "15.4 °C"

# Row 39
$ws.Cells.Item(39, 5).Value = "2026-02-11 23:20:05"
$ws.Cells.Item(39, 9).Value = "4.2 mm"

# Row 40
$ws.Cells.Item(40, 5).Value = "2026-02-11 23:20:07"
$ws.Cells.Item(40, 10).Value = "1006.1 hPa"

# Row 41
$ws.Cells.Item(41, 5).Value = "2026-02-11 23:20:10"
$ws.Cells.Item(41, 10).Value = "1004.0 hPa"

# Row 42
$ws.Cells.Item(42, 5).Value = "2026-02-11 23:20:12"

# Row 43
$ws.Cells.Item(43, 5).Value = "2026-02-11 23:20:15"
$ws.Cells.Item(43, 8).Value = "'69%"
$ws.Cells.Item(43, 15).Value = "12.6 °C"

# Row 44
$ws.Cells.Item(44, 5).Value = "2026-02-11 23:20:17"

# Row 45
$ws.Cells.Item(45, 5).Value = "2026-02-11 23:20:20"
$ws.Cells.Item(45, 10).Value = "1004.7 hPa"

# Row 46
$ws.Cells.Item(46, 5).Value = "2026-02-11 23:20:22"
$ws.Cells.Item(46, 10).Value = "1006.4 hPa"
$ws.Cells.Item(46, 12).Value = "48.2 km/h - 284º 22:51 TU"
